# Sync remote pairing matrix (rows 15-22) with newly logged pairing minutes,
# add a trailing marker cell, and move the active selection.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- Row 15 (Yuma Buchrieser) ---
$ws.Range("E15").Value = 660
$ws.Range("F15").Value = 150
$ws.Range("H15").Value = 90

# --- Row 16 (Wei Yi Stanley Ho) ---
$ws.Range("C16").Value = 60
$ws.Range("G16").Value = 660
$ws.Range("H16").Value = 660

# --- Row 17 (Dominik König) ---
$ws.Range("E17").Value = 150
$ws.Range("G17").Value = 330
$ws.Range("H17").Value = 450

# --- Row 18 (Martin Rabensteiner) ---
$ws.Range("F18").Value = 600

# --- Row 19 (Hannes Großauer) ---
$ws.Range("I19").Value = 660

# --- Row 20 (Rafael Forsthofer) ---
$ws.Range("I20").Value = 150

# Re-apply the SUM row as one fill operation so Excel stores it as a shared
# formula (matches how a drag-fill / paste of the row would serialize).
$ws.Range("B23:I23").Formula = "=SUM(B15:B22)"

# Trailing marker cell used as a layout placeholder further down the sheet.
$ws.Range("N30").Value = "                   "

# Move viewport / selection to match the latest edit location.
$ws.Application.ActiveWindow.ScrollRow = 10
$ws.Range("K32").Select()

$wb.Application.Calculate()
